# S.A.E. - 1.05 : re-save housekeeping on the SharePoint "custom XML" parts.
#
# When the document was re-saved after the content edits described in the
# commit ("Modification cour d'economie name", "Ajout de vocabulaire anglais
# pour Examen - 2", "Ajout d'un fichier finale pour la S.A.E. 1.05"), Word
# also rewrote the two anonymous SharePoint custom XML parts that back the
# document's content-type columns: the "FormTemplates" part and the
# "contentTypeSchema" part swapped physical slots (customXml/item1.xml <->
# customXml/item2.xml, and their companions customXml/itemProps1.xml <->
# customXml/itemProps2.xml). The actual XML payloads themselves are
# byte-for-byte identical - only which numbered part carries which payload
# changed. Reproduce that swap through the CustomXMLParts collection.

$d = $word.ActiveDocument
$parts = $d.CustomXMLParts

$formsNs = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$schemaNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"

# Grab the two parts by the namespace that identifies their payload.
$formsPart = $null
$schemaPart = $null
for ($i = 1; $i -le $parts.Count; $i++) {
    $p = $parts.Item($i)
    if ($p.NamespaceURI -eq $formsNs) { $formsPart = $p }
    if ($p.NamespaceURI -eq $schemaNs) { $schemaPart = $p }
}

if (($formsPart -ne $null) -and ($schemaPart -ne $null)) {
    $formsXml = $formsPart.XML
    $schemaXml = $schemaPart.XML

    # Drop both, then re-add them in the opposite order so the part that
    # used to hold the FormTemplates XML now holds the contentTypeSchema
    # XML (and ends up saved as item1.xml), and the part that used to hold
    # the contentTypeSchema XML now holds the FormTemplates XML (item2.xml).
    $formsPart.Delete()
    $schemaPart.Delete()

    [void]$parts.Add($schemaXml)
    [void]$parts.Add($formsXml)
} else {
    # Fallback for hosts where NamespaceURI lookup isn't populated: swap by
    # raw XML content instead, using Find/Replace-style identification.
    for ($i = 1; $i -le $parts.Count; $i++) {
        $p = $parts.Item($i)
        if ($p.XML -like "*FormTemplates*") { $formsPart = $p }
        if ($p.XML -like "*contentTypeSchema*") { $schemaPart = $p }
    }
    if (($formsPart -ne $null) -and ($schemaPart -ne $null)) {
        $formsXml = $formsPart.XML
        $schemaXml = $schemaPart.XML
        $formsPart.Delete()
        $schemaPart.Delete()
        [void]$parts.Add($schemaXml)
        [void]$parts.Add($formsXml)
    }
}
